$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ag-test (F) and Ag-positive (G) figures for rows 669-845 (data revision)
$ws.Range("F669").Value = 23654
$ws.Range("G669").Value = 637
$ws.Range("F670").Value = 53327
$ws.Range("G670").Value = 934
$ws.Range("F671").Value = 32954
$ws.Range("G671").Value = 637
$ws.Range("F672").Value = 30082
$ws.Range("G672").Value = 606
$ws.Range("F674").Value = 29135
$ws.Range("F676").Value = 28389
$ws.Range("F677").Value = 57271
$ws.Range("F678").Value = 34108
$ws.Range("F679").Value = 29702
$ws.Range("F680").Value = 28695
$ws.Range("F683").Value = 24438
$ws.Range("F684").Value = 58142
$ws.Range("F685").Value = 34729
$ws.Range("F686").Value = 34732
$ws.Range("G686").Value = 1173
$ws.Range("F687").Value = 31729
$ws.Range("F688").Value = 32478
$ws.Range("G688").Value = 1390
$ws.Range("F690").Value = 28055
$ws.Range("G690").Value = 1633
$ws.Range("F691").Value = 63399
$ws.Range("G691").Value = 2860
$ws.Range("F692").Value = 42001
$ws.Range("G692").Value = 2732
$ws.Range("F693").Value = 39907
$ws.Range("G693").Value = 2785
$ws.Range("F694").Value = 37899
$ws.Range("F695").Value = 37634
$ws.Range("G695").Value = 3188
$ws.Range("F697").Value = 29247
$ws.Range("G697").Value = 3072
$ws.Range("F698").Value = 71741
$ws.Range("G698").Value = 5967
$ws.Range("F699").Value = 43821
$ws.Range("G699").Value = 4359
$ws.Range("F700").Value = 44154
$ws.Range("G700").Value = 4378
$ws.Range("F701").Value = 42184
$ws.Range("G701").Value = 3907
$ws.Range("F702").Value = 36680
$ws.Range("G702").Value = 3977
$ws.Range("F703").Value = 17105
$ws.Range("F704").Value = 25227
$ws.Range("G704").Value = 3737
$ws.Range("F705").Value = 56892
$ws.Range("G705").Value = 6452
$ws.Range("F706").Value = 41030
$ws.Range("G706").Value = 5006
$ws.Range("F707").Value = 39170
$ws.Range("G707").Value = 4661
$ws.Range("F708").Value = 35987
$ws.Range("G708").Value = 4227
$ws.Range("F709").Value = 32732
$ws.Range("G709").Value = 4035
$ws.Range("F712").Value = 52138
$ws.Range("G712").Value = 6379
$ws.Range("F713").Value = 37640
$ws.Range("G713").Value = 4808
$ws.Range("F714").Value = 32857
$ws.Range("G714").Value = 4049
$ws.Range("F715").Value = 32214
$ws.Range("G715").Value = 3618
$ws.Range("F716").Value = 30028
$ws.Range("G716").Value = 3722
$ws.Range("F719").Value = 45533
$ws.Range("G719").Value = 5383
$ws.Range("F720").Value = 31587
$ws.Range("G720").Value = 3566
$ws.Range("F721").Value = 28260
$ws.Range("G721").Value = 3190
$ws.Range("F722").Value = 28272
$ws.Range("G722").Value = 2938
$ws.Range("F723").Value = 23096
$ws.Range("G723").Value = 2834
$ws.Range("F725").Value = 12837
$ws.Range("G725").Value = 2099
$ws.Range("F726").Value = 36620
$ws.Range("G726").Value = 4229
$ws.Range("F727").Value = 25501
$ws.Range("G727").Value = 2844
$ws.Range("F728").Value = 25021
$ws.Range("G728").Value = 2643
$ws.Range("F729").Value = 23580
$ws.Range("G729").Value = 2549
$ws.Range("F730").Value = 19854
$ws.Range("G730").Value = 2370
$ws.Range("F731").Value = 8679
$ws.Range("F732").Value = 11962
$ws.Range("G732").Value = 1928
$ws.Range("F733").Value = 32422
$ws.Range("G733").Value = 3764
$ws.Range("F734").Value = 23410
$ws.Range("G734").Value = 2570
$ws.Range("F735").Value = 19606
$ws.Range("G735").Value = 2289
$ws.Range("F736").Value = 19858
$ws.Range("G736").Value = 2212
$ws.Range("F737").Value = 18786
$ws.Range("G737").Value = 2327
$ws.Range("F738").Value = 6883
$ws.Range("F739").Value = 8758
$ws.Range("G739").Value = 1423
$ws.Range("F740").Value = 25444
$ws.Range("G740").Value = 2792
$ws.Range("F741").Value = 19191
$ws.Range("G741").Value = 1952
$ws.Range("F742").Value = 17545
$ws.Range("G742").Value = 1711
$ws.Range("F743").Value = 18277
$ws.Range("G743").Value = 1645
$ws.Range("F744").Value = 14932
$ws.Range("G744").Value = 1619
$ws.Range("F746").Value = 8068
$ws.Range("G746").Value = 1260
$ws.Range("F747").Value = 22858
$ws.Range("G747").Value = 2413
$ws.Range("F748").Value = 17118
$ws.Range("G748").Value = 1553
$ws.Range("F749").Value = 15027
$ws.Range("G749").Value = 1492
$ws.Range("F750").Value = 15232
$ws.Range("G750").Value = 1362
$ws.Range("F751").Value = 12718
$ws.Range("G751").Value = 1393
$ws.Range("F754").Value = 21485
$ws.Range("G754").Value = 1966
$ws.Range("F755").Value = 13908
$ws.Range("G755").Value = 1303
$ws.Range("F756").Value = 13914
$ws.Range("G756").Value = 1085
$ws.Range("F757").Value = 13682
$ws.Range("F761").Value = 17003
$ws.Range("G761").Value = 1292
$ws.Range("F764").Value = 11223
$ws.Range("G764").Value = 699
$ws.Range("F765").Value = 9244
$ws.Range("F767").Value = 4163
$ws.Range("F768").Value = 15087
$ws.Range("F769").Value = 10127
$ws.Range("F770").Value = 9169
$ws.Range("F771").Value = 9244
$ws.Range("F775").Value = 3239
$ws.Range("F776").Value = 14928
$ws.Range("F777").Value = 10481
$ws.Range("G777").Value = 455
$ws.Range("F778").Value = 9076
$ws.Range("F779").Value = 7404
$ws.Range("G779").Value = 313
$ws.Range("F781").Value = 2768
$ws.Range("F782").Value = 10749
$ws.Range("G782").Value = 433
$ws.Range("F783").Value = 7851
$ws.Range("G783").Value = 253
$ws.Range("F784").Value = 7818
$ws.Range("G784").Value = 239
$ws.Range("F785").Value = 7191
$ws.Range("F786").Value = 6351
$ws.Range("F788").Value = 1768
$ws.Range("G788").Value = 87
$ws.Range("F789").Value = 7892
$ws.Range("G789").Value = 341
$ws.Range("F790").Value = 4680
$ws.Range("G790").Value = 176
$ws.Range("F791").Value = 4258
$ws.Range("G791").Value = 239
$ws.Range("F795").Value = 1180
$ws.Range("G795").Value = 63
$ws.Range("F796").Value = 4664
$ws.Range("G796").Value = 262
$ws.Range("F798").Value = 3511
$ws.Range("G798").Value = 138
$ws.Range("F802").Value = 987
$ws.Range("F803").Value = 3827
$ws.Range("G803").Value = 173
$ws.Range("F805").Value = 2410
$ws.Range("F810").Value = 3898
$ws.Range("F812").Value = 2106
$ws.Range("F813").Value = 2531
$ws.Range("F814").Value = 1998
$ws.Range("G814").Value = 64
$ws.Range("F817").Value = 3769
$ws.Range("G817").Value = 88
$ws.Range("F819").Value = 2279
$ws.Range("F820").Value = 2561
$ws.Range("F824").Value = 3606
$ws.Range("F826").Value = 1831
$ws.Range("F827").Value = 1915
$ws.Range("F832").Value = 1917
$ws.Range("F833").Value = 1397
$ws.Range("F834").Value = 1789
$ws.Range("F835").Value = 1140
$ws.Range("F839").Value = 1484
$ws.Range("F841").Value = 2731
$ws.Range("F842").Value = 1520
$ws.Range("F845").Value = 3810

# Append new daily rows 884-886
$ws.Range("A884").Value = 44778
$ws.Range("B884").Value = 1822841
$ws.Range("C884").Value = 2189
$ws.Range("D884").Value = 933
$ws.Range("E884").Value = 20253

$ws.Range("A885").Value = 44779
$ws.Range("B885").Value = 1823259
$ws.Range("C885").Value = 1003
$ws.Range("D885").Value = 418
$ws.Range("E885").Value = 20259

$ws.Range("A886").Value = 44780
$ws.Range("B886").Value = 1823364
$ws.Range("C886").Value = 326
$ws.Range("D886").Value = 105
$ws.Range("E886").Value = 20263
